$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "IdentityDocumentTypes Data")
$ws.Name = "Data"

# Add the new "Options" column header
$ws.Range("C1").Value = "Options"

# Fill in the "Options" value for each existing identity document type row
$ws.Range("C2").Value = "IsRussianBirthCertificate"
$ws.Range("C3").Value = "IsRussianPassport"
$ws.Range("C4").Value = "NULL"
$ws.Range("C5").Value = "NULL"

# Append the new identity document type row.
# Force the Id to be stored as text (matching the existing "1".."4" text
# values in column A) instead of being auto-coerced to a number.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "5"
$ws.Range("A6").ClearFormats()

$ws.Range("B6").Value = "Заграничный паспорт гражданина РФ"
$ws.Range("C6").Value = "IsRussianForeignPassport"
